# adding goals info back to the df and datetime FE
#
# Insert 4 new rows (homeGoals, awayGoals, homeGoalsHalfTime,
# awayGoalsHalfTime) right after the existing awayTeamID row (row 7),
# pushing all the existing stat rows down by 4. Final row count grows
# from 44 to 48 (dimension A1:B48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the old row 8 ("home_xGoals"), shifting
# everything from row 8 down to row 12 onward.
$ws.Rows("8:11").Insert()

# Carry over the header cell-style (bold/bordered/centered label style)
# from the row directly above onto the newly inserted label cells, so
# the new rows look like the rest of the key/value table.
$ws.Range("A7").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)

# Populate the newly inserted rows with the goals info.
$ws.Range("A8").Value = "homeGoals"
$ws.Range("B8").Value = 11

$ws.Range("A9").Value = "awayGoals"
$ws.Range("B9").Value = 10

$ws.Range("A10").Value = "homeGoalsHalfTime"
$ws.Range("B10").Value = 7

$ws.Range("A11").Value = "awayGoalsHalfTime"
$ws.Range("B11").Value = 6
